$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the same 40-character column width used by A:D to the new
# reporting columns E:H so the whole table is formatted consistently.
$ws.Columns("E:H").ColumnWidth = 40

# Row 1 - headers (names)
$ws.Range("A1").Value = "arango juan"
$ws.Range("B1").Value = "Gary Tsai"
$ws.Range("C1").Value = "Covalky Pena"
$ws.Range("D1").Value = "Crystal Evelyn "
$ws.Range("E1").Value = "Jean Augustin"
$ws.Range("F1").Value = "Natalie Primus"
$ws.Range("G1").Value = "Miguel Martillo"
$ws.Range("H1").Value = "Anne Crosby"

# Row 2 - codes
$ws.Range("A2").Value = "CFD8AFA4C0"
$ws.Range("B2").Value = "CFD893A460"
$ws.Range("C2").Value = "0FD8AD42A0"
$ws.Range("D2").Value = "0FD8A87380"
$ws.Range("E2").Value = "4FD8B36A40"
$ws.Range("F2").Value = "0FD8AE8B60"
$ws.Range("G2").Value = "8FD8B68DE0"
$ws.Range("H2").Value = "4FD8A33DE0"

# Row 3 - IN times
$ws.Range("A3").Value = "IN -> 2017/01/31 18:57"
$ws.Range("B3").Value = "IN -> 2017/01/31 19:12"
$ws.Range("C3").Value = "IN -> 2017/01/31 19:19"
$ws.Range("D3").Value = "IN -> 2017/01/31 19:23"
$ws.Range("E3").Value = "IN -> 2017/01/31 19:24"
$ws.Range("F3").Value = "IN -> 2017/02/01 08:49"
$ws.Range("G3").Value = "IN -> 2017/02/01 09:06"
$ws.Range("H3").Value = "IN -> 2017/02/01 09:08"

# Row 4 - OUT times
$ws.Range("A4").Value = "OUT -> 2017/01/31 19:07"
$ws.Range("B4").Value = "OUT -> 2017/01/31 19:12"
$ws.Range("C4").Value = "OUT -> 2017/02/01 15:25"

# Row 5
$ws.Range("A5").Value = "IN -> 2017/02/01 17:09"
$ws.Range("B5").Value = "IN -> 2017/01/31 19:19"
$ws.Range("C5").Value = "IN -> 2017/02/01 17:52"

# Row 6-10 (only column B)
$ws.Range("B6").Value = "OUT -> 2017/01/31 19:23"
$ws.Range("B7").Value = "IN -> 2017/01/31 19:26"
$ws.Range("B8").Value = "OUT -> 2017/01/31 19:33"
$ws.Range("B9").Value = "IN -> 2017/02/01 15:01"
$ws.Range("B10").Value = "OUT -> 2017/02/01 17:31"
